# Insert one new price-report row for "Acelga" (Vega Modelo de Temuco) at
# row 354, pushing the existing rows 354-440 down to 355-441 (dimension
# grows from A1:R440 to A1:R441).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 354 (shifts everything below it down by one).
$ws.Rows.Item(354).Insert()

# Populate the newly inserted row 354 with the new record.
$ws.Cells.Item(354, 1).Value  = 10
$ws.Cells.Item(354, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(354, 3).Value  = "La Araucanía"
$ws.Cells.Item(354, 4).Value  = 44943
$ws.Cells.Item(354, 5).Value  = 9
$ws.Cells.Item(354, 6).Value  = 100112009
$ws.Cells.Item(354, 7).Value  = "Acelga"
$ws.Cells.Item(354, 8).Value  = "Sin especificar"
$ws.Cells.Item(354, 9).Value  = "Primera"
$ws.Cells.Item(354, 10).Value = 120
$ws.Cells.Item(354, 11).Value = 8000
$ws.Cells.Item(354, 12).Value = 9000
$ws.Cells.Item(354, 13).Value = 8458
$ws.Cells.Item(354, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(354, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(354, 16).Value = 705
$ws.Cells.Item(354, 17).Value = 12
$ws.Cells.Item(354, 18).Value = "Hortaliza"
